$d = $word.ActiveDocument

# 1. Remove the "Meta description: ..." paragraph (originally paragraph 2,
#    right after the title heading).
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# 2. Insert a new bold paragraph "Play Chicago Gold Free: Thrilling Online
#    Slot game" right before the final paragraph (the one that used to hold
#    the "Prompt: ..." image-generation text).
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$lastPara.Range.InsertParagraphBefore()

$newCount = $d.Paragraphs.Count
$newPara = $d.Paragraphs.Item($newCount - 1)
$npr = $newPara.Range
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Chicago Gold Free: Thrilling Online Slot game</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$npr.InsertXML($xml)

# 3. Replace the old image-prompt text in the final paragraph with the new
#    review blurb, keeping the existing italic run formatting.
$old = 'Prompt: Create a feature image fitting the game "Chicago Gold". The image should be in cartoon style and should feature a happy Maya warrior with glasses. The Maya warrior should be portrayed as a successful thief, surrounded by diamonds and holding a bag of stolen goods. The background of the image should be a night scene of the city of Chicago with a spotlight shining on the Maya warrior. The overall tone of the image should be fun and adventurous, representing the excitement of playing the game.'
$new = 'Read our review of Chicago Gold, the thrilling online slot game set in the criminal underworld of Chicago. Play now for free and escape with riches.'
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
